$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.487.12'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.37%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.809.00'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.04%  '
$ws.Range("E4").Value = '  -0.41%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.004'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '307.14'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.21%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4525'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.16%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3602'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.88%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.54'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.48%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07101'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.91%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8913'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07810'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.00%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.52'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.11%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.831.87'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.23%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.299'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.37%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.339'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.20%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '85.25'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.16%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.006'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.31%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008489'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.003'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '26.518.76'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.34%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.28'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.977'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.038.39'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.53'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.81%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.968'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '151.15'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.29%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.83'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.65%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.058'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '112.13'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.35%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.880'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.56%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08703'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.82%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.120'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.51%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.869'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +14.87%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.446'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.58%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7241'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.07%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.114'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.24%  '
$ws.Range("E38").Value = '  -0.31%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.073'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.86%  '
$ws.Range("E40").Value = '  +0.23%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.05108'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.23%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.897'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.52%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5159'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.804'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.47%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1515'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.66%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.037'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.25%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4671'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.62%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.003'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.38%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.976'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.55%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '100.61'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.28%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.576'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.07%  '
